# Update "想去人数" (F column) values across the workbook's sheets.
$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 46
$ws1.Range("F4").Value  = 643
$ws1.Range("F5").Value  = 493
$ws1.Range("F6").Value  = 1164
$ws1.Range("F8").Value  = 183
$ws1.Range("F9").Value  = 68
$ws1.Range("F10").Value = 770
$ws1.Range("F11").Value = 418
$ws1.Range("F15").Value = 205
$ws1.Range("F16").Value = 9
$ws1.Range("F17").Value = 387
$ws1.Range("F18").Value = 6274
$ws1.Range("F20").Value = 58
$ws1.Range("F21").Value = 15
$ws1.Range("F22").Value = 7240
$ws1.Range("F25").Value = 3298
$ws1.Range("F26").Value = 428
$ws1.Range("F27").Value = 804
$ws1.Range("F28").Value = 4479
$ws1.Range("F29").Value = 334
$ws1.Range("F30").Value = 165
$ws1.Range("F31").Value = 161
$ws1.Range("F32").Value = 1304
$ws1.Range("F36").Value = 1021
$ws1.Range("F37").Value = 1337
$ws1.Range("F38").Value = 2084

# --- Sheet "演出" ---
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 36

# --- Sheet "本地生活" ---
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F3").Value = 1172
$ws3.Range("F4").Value = 63

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 46
$ws4.Range("F4").Value  = 1172
$ws4.Range("F5").Value  = 63
$ws4.Range("F7").Value  = 643
$ws4.Range("F8").Value  = 493
$ws4.Range("F9").Value  = 1164
$ws4.Range("F11").Value = 183
$ws4.Range("F12").Value = 68
$ws4.Range("F13").Value = 770
$ws4.Range("F14").Value = 418
$ws4.Range("F19").Value = 205
$ws4.Range("F20").Value = 9
$ws4.Range("F21").Value = 387
$ws4.Range("F22").Value = 6274
$ws4.Range("F24").Value = 58
$ws4.Range("F25").Value = 15
$ws4.Range("F26").Value = 7240
$ws4.Range("F29").Value = 3298
$ws4.Range("F30").Value = 428
$ws4.Range("F31").Value = 804
$ws4.Range("F32").Value = 4479
$ws4.Range("F33").Value = 334
$ws4.Range("F34").Value = 36
$ws4.Range("F35").Value = 165
$ws4.Range("F36").Value = 161
$ws4.Range("F37").Value = 1304
$ws4.Range("F41").Value = 1021
$ws4.Range("F42").Value = 1337
$ws4.Range("F44").Value = 2084
